$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 543
$ws.Range("I2").Value = 1340
$ws.Range("J2").Value = 5817
$ws.Range("K2").Value = 31
$ws.Range("M2").Value = 91
$ws.Range("N2").Value = 980
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 24
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 82
$ws.Range("S2").Value = 601
$ws.Range("T2").Value = 933
$ws.Range("U2").Value = 68
$ws.Range("V2").Value = 8992
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 8816
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 141
$ws.Range("AA2").Value = 52
